# Apply "last minute updates" to the first paragraph of the document:
#   - add a paragraph border (5pt space on all four sides)
#   - increase the left indent from 6pt (120 twips) to 11.25pt (225 twips)
#   - drop the trailing " " run
#   - rename the placeholder id from ..._topic_3__ID** to ..._104__ID**

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$pFmt = $p1.Range.ParagraphFormat

# Paragraph border: <w:pBdr><w:top w:space="5"/><w:left w:space="5"/>
#                    <w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$borders = $pFmt.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# <w:ind w:left="225"/> (was 120 => 6pt, now 225 => 11.25pt)
$pFmt.LeftIndent = 11.25

# Remove the trailing space-only run, leaving a single run in the paragraph.
$parEnd = $p1.Range.End
$trailingSpace = $d.Range($parEnd - 2, $parEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# Rename the placeholder id text.
$d.Content.Find.Execute("**ID__AFFARS_5337_topic_3__ID**", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "**ID__AFFARS_5337_104__ID**", 2)
